# "Generate Report for Handback" - mark the zh-cn and de-de localization
# rows as handed back, fill in the target/handback file links + datetime
# stamps, and widen the columns that now hold longer content.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

$srcMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/77d063ccd3225f60e744edc0d3d8702ce9c78293/e2e/1615029a-e0f8-439e-b6b7-f7b1579e7152.md"
$srcMdName = "1615029a-e0f8-439e-b6b7-f7b1579e7152.md"
$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: status shown per-language ----
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Columns.Item(5).ColumnWidth = 29.1666666666667
$ws1.Columns.Item(6).ColumnWidth = 29.1666666666667

# ---- zh-cn sheet ----
$ws2.Range("C2").Value = $newStatus

$ws2.Hyperlinks.Add($ws2.Range("I2"), $srcMdUrl, "", "", $srcMdName)
$ws2.Cells.Item(2, 9).Font.Name = "Calibri"
$ws2.Cells.Item(2, 9).Font.Color = 15570276

$ws2.Range("J2").Value = "1615029a-e0f8-439e-b6b7-f7b1579e7152.2783e4b09fd901c8a7f5cac97f593e10abcd79d8.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-07 14:35:02"

$ws2.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws2.Columns.Item(9).ColumnWidth = 39.1666666666667
$ws2.Columns.Item(10).ColumnWidth = 39.1666666666667

# ---- de-de sheet ----
$ws3.Range("C2").Value = $newStatus

$ws3.Hyperlinks.Add($ws3.Range("I2"), $srcMdUrl, "", "", $srcMdName)
$ws3.Cells.Item(2, 9).Font.Name = "Calibri"
$ws3.Cells.Item(2, 9).Font.Color = 15570276

$ws3.Range("J2").Value = "1615029a-e0f8-439e-b6b7-f7b1579e7152.2783e4b09fd901c8a7f5cac97f593e10abcd79d8.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-07 14:35:48"

$ws3.Columns.Item(3).ColumnWidth = 29.1666666666667
$ws3.Columns.Item(9).ColumnWidth = 39.1666666666667
$ws3.Columns.Item(10).ColumnWidth = 39.1666666666667
